$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.211.15"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.828.82"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'234.22"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "'0.6000"
$ws.Range("E6").Value = "  -4.31%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "'0.06979"
$ws.Range("E8").Value = "  -5.76%  "
$ws.Range("D9").Value = "'0.2757"
$ws.Range("E9").Value = "  -4.74%  "
$ws.Range("D10").Value = "'23.24"
$ws.Range("E10").Value = "  -6.76%  "
$ws.Range("D11").Value = "'0.07614"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "1.828.89"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "'4.755"
$ws.Range("E13").Value = "  -4.22%  "
$ws.Range("D14").Value = "'0.6271"
$ws.Range("E14").Value = "  -6.40%  "
$ws.Range("D15").Value = "'0.000009673"
$ws.Range("E15").Value = "  -6.97%  "
$ws.Range("D16").Value = "'78.33"
$ws.Range("E16").Value = "  -4.27%  "
$ws.Range("D17").Value = "28.836.18"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "'5.708"
$ws.Range("E18").Value = "  -8.75%  "
$ws.Range("D19").Value = "'220.63"
$ws.Range("E19").Value = "  -5.76%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("D21").Value = "'11.55"
$ws.Range("E21").Value = "  -6.08%  "
$ws.Range("D22").Value = "'6.847"
$ws.Range("E22").Value = "  -6.24%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'155.75"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").Value = "'7.961"
$ws.Range("E25").Value = "  -6.05%  "
$ws.Range("D26").Value = "'0.1291"
$ws.Range("E26").Value = "  -4.22%  "
$ws.Range("D27").Value = "'16.53"
$ws.Range("E27").Value = "  -4.54%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'1.449"
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.06403"
$ws.Range("E29").Value = "  -11.96%  "
$ws.Range("D30").Value = "'1.440"
$ws.Range("E30").Value = "  -2.90%  "
$ws.Range("D31").Value = "'3.833"
$ws.Range("E31").Value = "  -4.98%  "
$ws.Range("D32").Value = "'3.762"
$ws.Range("E32").Value = "  -6.92%  "
$ws.Range("D33").Value = "'1.092"
$ws.Range("E33").Value = "  -5.95%  "
$ws.Range("D34").Value = "'1.722"
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("D35").Value = "'0.6444"
$ws.Range("E35").Value = "  -9.93%  "
$ws.Range("D36").Value = "'2.534"
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("D37").Value = "'2.732"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "'0.01748"
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("D39").Value = "'6.538"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("D40").Value = "1.171.42"
$ws.Range("E40").Value = "  -4.93%  "
$ws.Range("D41").Value = "'0.8937"
$ws.Range("E41").Value = "  -6.60%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "1.983.14"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'100.70"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "'62.12"
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("E46").Value = "  -5.07%  "
$ws.Range("D47").Value = "'0.05639"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "'8.459"
$ws.Range("E48").Value = "  -5.03%  "
$ws.Range("D49").Value = "'1.585"
$ws.Range("E49").Value = "  -6.78%  "
$ws.Range("D50").Value = "'0.4551"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").Value = "'0.3647"
$ws.Range("E51").Value = "  -6.11%  "
